# added lnu 2020 data
# The original sheet had two rows that needed to be removed:
#  - row 12, a duplicate entry (same DOI as row 11, wrong amount)
#  - row 76, an extra entry that should not have been included
# Deleting them shifts every following row up by one (or two), which is
# exactly what the target workbook looks like.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the earlier deletion's row
# index (12) is not affected by this one shifting rows upward.
$ws.Rows.Item(76).Delete()
$ws.Rows.Item(12).Delete()

# Restore the view state recorded in the saved workbook (scrolled down a
# bit, with a single active cell selected instead of a range).
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("C43").Select()
